$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1744.5333
$ws.Range("I92").Value = 1276.625
$ws.Range("J92").Value = 2279.2856
$ws.Range("K92").Value = 1276.625
$ws.Range("L92").Value = 2279.2856
$ws.Range("M92").Value = -28.625
$ws.Range("N92").Value = -4775.2856
$ws.Range("H94").Value = 45566596
$ws.Range("I94").Value = 62500320
$ws.Range("K94").Value = 62500320
$ws.Range("M94").Value = -62499869
$ws.Range("H125").Value = 3884.85
$ws.Range("I125").Value = 5124.3335
$ws.Range("K125").Value = 46119.0015
$ws.Range("M125").Value = -43659.0015
$ws.Range("H135").Value = 2238.4375
$ws.Range("I135").Value = 2308.1785
$ws.Range("K135").Value = 20773.6065
$ws.Range("M135").Value = -18238.6065
$ws.Range("H138").Value = 1794.8306
$ws.Range("I138").Value = 1072.1395
$ws.Range("J138").Value = 3737.0625
$ws.Range("K138").Value = 3216.4185
$ws.Range("L138").Value = 11211.1875
$ws.Range("M138").Value = 1923.5815
$ws.Range("N138").Value = -21491.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5693.794
$ws.Range("I32").Value = 5818.123
$ws.Range("K32").Value = 5818.123
$ws.Range("M32").Value = -5531.123
$ws.Range("H45").Value = 8923.556
$ws.Range("I45").Value = 8262.4
$ws.Range("K45").Value = 8262.4
$ws.Range("M45").Value = -7885.4
$ws.Range("H61").Value = 4790.3857
$ws.Range("I61").Value = 4803.9243
$ws.Range("K61").Value = 4803.9243
$ws.Range("M61").Value = -4591.9243
$ws.Range("H74").Value = 2318.1667
$ws.Range("I74").Value = 893.94446
$ws.Range("J74").Value = 6590.8335
$ws.Range("K74").Value = 893.94446
$ws.Range("L74").Value = 6590.8335
$ws.Range("M74").Value = -19.94446000000005
$ws.Range("N74").Value = -8338.833500000001
$ws.Range("H77").Value = 2318.1667
$ws.Range("I77").Value = 893.94446
$ws.Range("J77").Value = 6590.8335
$ws.Range("K77").Value = 4469.7223
$ws.Range("L77").Value = 32954.1675
$ws.Range("M77").Value = -101.7223000000004
$ws.Range("N77").Value = -41690.1675
$ws.Range("H97").Value = 5197.5654
$ws.Range("I97").Value = 6315.778
$ws.Range("K97").Value = 6315.778
$ws.Range("M97").Value = -5819.778
$ws.Range("H102").Value = 10493
$ws.Range("I102").Value = 13982.75
$ws.Range("K102").Value = 13982.75
$ws.Range("M102").Value = -12360.75
$ws.Range("H130").Value = 29999.334
$ws.Range("J130").Value = 29999.334
$ws.Range("L130").Value = 29999.334
$ws.Range("N130").Value = -40039.334
$ws.Range("H136").Value = 4790.3857
$ws.Range("I136").Value = 4803.9243
$ws.Range("K136").Value = 14411.7729
$ws.Range("M136").Value = -11861.7729

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 4319.8
$ws.Range("J8").Value = 4319.8
$ws.Range("L8").Value = 4319.8
$ws.Range("N8").Value = -4599.8
$ws.Range("H22").Value = 1021.5
$ws.Range("I22").Value = 1028.6666
$ws.Range("K22").Value = 1028.6666
$ws.Range("M22").Value = -855.6666
$ws.Range("H105").Value = 3894
$ws.Range("I105").Value = 1641.3
$ws.Range("J105").Value = 8399.4
$ws.Range("K105").Value = 1641.3
$ws.Range("L105").Value = 8399.4
$ws.Range("M105").Value = 105.7
$ws.Range("N105").Value = -11893.4
$ws.Range("H134").Value = 7040.1226
$ws.Range("I134").Value = 7181.0454
$ws.Range("K134").Value = 21543.1362
$ws.Range("M134").Value = -19008.1362
$ws.Range("H141").Value = 20000
$ws.Range("I141").Value = 20000
$ws.Range("K141").Value = 20000
$ws.Range("M141").Value = -14820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3610.3438
$ws.Range("I31").Value = 2530.2
$ws.Range("J31").Value = 4563.4116
$ws.Range("K31").Value = 2530.2
$ws.Range("L31").Value = 4563.4116
$ws.Range("M31").Value = -2235.2
$ws.Range("N31").Value = -5153.4116
$ws.Range("H34").Value = 3610.3438
$ws.Range("I34").Value = 2530.2
$ws.Range("J34").Value = 4563.4116
$ws.Range("K34").Value = 2530.2
$ws.Range("L34").Value = 4563.4116
$ws.Range("M34").Value = -2328.2
$ws.Range("N34").Value = -4967.4116
$ws.Range("H62").Value = 44494.438
$ws.Range("J62").Value = 80695.375
$ws.Range("L62").Value = 80695.375
$ws.Range("N62").Value = -81943.375
$ws.Range("H65").Value = 44494.438
$ws.Range("J65").Value = 80695.375
$ws.Range("L65").Value = 403476.875
$ws.Range("N65").Value = -409716.875
$ws.Range("H107").Value = 9019.925999999999
$ws.Range("I107").Value = 12432.947
$ws.Range("K107").Value = 12432.947
$ws.Range("M107").Value = -10512.947
$ws.Range("H122").Value = 19258.715
$ws.Range("I122").Value = 26202.4
$ws.Range("J122").Value = 1899.5
$ws.Range("K122").Value = 78607.20000000001
$ws.Range("L122").Value = 5698.5
$ws.Range("M122").Value = -76157.20000000001
$ws.Range("N122").Value = -10598.5
$ws.Range("H134").Value = 1653.0857
$ws.Range("J134").Value = 1898.8
$ws.Range("L134").Value = 5696.4
$ws.Range("N134").Value = -10766.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 66135228
$ws.Range("I4").Value = 56659396
$ws.Range("K4").Value = 169978188
$ws.Range("M4").Value = -169978076
$ws.Range("H5").Value = 527763.0600000001
$ws.Range("I5").Value = 583.6667
$ws.Range("J5").Value = 771076.6
$ws.Range("K5").Value = 1751.0001
$ws.Range("L5").Value = 2313229.8
$ws.Range("M5").Value = -1639.0001
$ws.Range("N5").Value = -2313453.8
$ws.Range("H44").Value = 1889.7778
$ws.Range("J44").Value = 5229
$ws.Range("L44").Value = 15687
$ws.Range("N44").Value = -16483
$ws.Range("H63").Value = 2993.5
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 2993.5
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H70").Value = 668
$ws.Range("I70").Value = 668
$ws.Range("K70").Value = 2004
$ws.Range("M70").Value = -1689
$ws.Range("H73").Value = 668
$ws.Range("I73").Value = 668
$ws.Range("K73").Value = 2004
$ws.Range("M73").Value = -912
$ws.Range("H107").Value = 945.3200000000001
$ws.Range("I107").Value = 308.33334
$ws.Range("J107").Value = 1146.4736
$ws.Range("K107").Value = 925.0000200000001
$ws.Range("L107").Value = 3439.4208
$ws.Range("M107").Value = 994.9999799999999
$ws.Range("N107").Value = -7279.4208
$ws.Range("H135").Value = 527763.0600000001
$ws.Range("I135").Value = 583.6667
$ws.Range("J135").Value = 771076.6
$ws.Range("K135").Value = 5253.0003
$ws.Range("L135").Value = 6939689.399999999
$ws.Range("M135").Value = -2718.0003
$ws.Range("N135").Value = -6944759.399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12467.056
$ws.Range("I122").Value = 7567
$ws.Range("K122").Value = 22701
$ws.Range("M122").Value = -20251
$ws.Range("H134").Value = 73264.625
$ws.Range("J134").Value = 73264.625
$ws.Range("L134").Value = 219793.875
$ws.Range("N134").Value = -224863.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19653.285
$ws.Range("I7").Value = 34081
$ws.Range("K7").Value = 34081
$ws.Range("M7").Value = -33969
$ws.Range("H40").Value = 17189
$ws.Range("I40").Value = 19518.092
$ws.Range("J40").Value = 12065
$ws.Range("K40").Value = 19518.092
$ws.Range("L40").Value = 12065
$ws.Range("M40").Value = -19382.092
$ws.Range("N40").Value = -12337
$ws.Range("H46").Value = 2500.3333
$ws.Range("I46").Value = 366
$ws.Range("J46").Value = 3567.5
$ws.Range("K46").Value = 366
$ws.Range("L46").Value = 3567.5
$ws.Range("M46").Value = -178
$ws.Range("N46").Value = -3943.5
$ws.Range("H126").Value = 19653.285
$ws.Range("I126").Value = 34081
$ws.Range("K126").Value = 102243
$ws.Range("M126").Value = -99773
$ws.Range("H132").Value = 518116.9
$ws.Range("I132").Value = 879090.0600000001
$ws.Range("J132").Value = 6738.25
$ws.Range("K132").Value = 2637270.18
$ws.Range("L132").Value = 20214.75
$ws.Range("M132").Value = -2634740.18
$ws.Range("N132").Value = -25274.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H100").Value = 33713.91
$ws.Range("I100").Value = 24586.467
$ws.Range("J100").Value = 53272.715
$ws.Range("K100").Value = 49172.934
$ws.Range("L100").Value = 106545.43
$ws.Range("M100").Value = -48631.934
$ws.Range("N100").Value = -107627.43
$ws.Range("H122").Value = 15987.073
$ws.Range("I122").Value = 2088.25
$ws.Range("K122").Value = 6264.75
$ws.Range("M122").Value = -3814.75
$ws.Range("H132").Value = 5477
$ws.Range("I132").Value = 6267.4717
$ws.Range("K132").Value = 18802.4151
$ws.Range("M132").Value = -16272.4151
$ws.Range("H136").Value = 359481.44
$ws.Range("I136").Value = 497537.4
$ws.Range("J136").Value = 2836.8333
$ws.Range("K136").Value = 1492612.2
$ws.Range("L136").Value = 8510.499899999999
$ws.Range("M136").Value = -1490062.2
$ws.Range("N136").Value = -13610.4999

Write-Host "done"